$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.392285
$ws.Range("H2").Value = 1.176855
$ws.Range("I2").Value = 0.0216687106565248
$ws.Range("J2").Value = 0.0216687106565248
$ws.Range("M2").Value = 72.266001
$ws.Range("N2").Value = 216.798003
$ws.Range("O2").Value = 0.2949652269937106
$ws.Range("P2").Value = 0.2949652269937106
$ws.Range("Q2").Value = 28.348868202285
$ws.Range("R2").Value = 255.139813820565
$ws.Range("S2").Value = 0.006391516157462874
$ws.Range("T2").Value = 0.006391516157462873

$ws.Range("G3").Value = 0.392285
$ws.Range("H3").Value = 1.176855
$ws.Range("I3").Value = 0.0216687106565248
$ws.Range("J3").Value = 0.0216687106565248
$ws.Range("N3").Value = 410.023338
$ws.Range("O3").Value = 0.5578585839920717
$ws.Range("P3").Value = 0.5578585839920718
$ws.Range("Q3").Value = 53.61533504910999
$ws.Range("R3").Value = 482.53801544199
$ws.Range("S3").Value = 0.01208807624378284
$ws.Range("T3").Value = 0.01208807624378284

$ws.Range("G4").Value = 0.392285
$ws.Range("H4").Value = 1.176855
$ws.Range("I4").Value = 0.0216687106565248
$ws.Range("J4").Value = 0.0216687106565248
$ws.Range("M4").Value = 36.057927
$ws.Range("N4").Value = 108.173781
$ws.Range("O4").Value = 0.1471761890142177
$ws.Range("P4").Value = 0.1471761890142177
$ws.Range("Q4").Value = 14.144983893195
$ws.Range("R4").Value = 127.304855038755
$ws.Range("S4").Value = 0.003189118255279088
$ws.Range("T4").Value = 0.003189118255279087

$ws.Range("I5").Value = 0.7688940604785444
$ws.Range("J5").Value = 0.7688940604785441
$ws.Range("M5").Value = 72.266001
$ws.Range("N5").Value = 216.798003
$ws.Range("O5").Value = 0.2949652269937106
$ws.Range("P5").Value = 0.2949652269937106
$ws.Range("Q5").Value = 1005.933242985202
$ws.Range("R5").Value = 9053.399186866818
$ws.Range("S5").Value = 0.2267970110831697
$ws.Range("T5").Value = 0.2267970110831696

$ws.Range("I6").Value = 0.7688940604785444
$ws.Range("J6").Value = 0.7688940604785441
$ws.Range("N6").Value = 410.023338
$ws.Range("O6").Value = 0.5578585839920717
$ws.Range("P6").Value = 0.5578585839920718
$ws.Range("Q6").Value = 1902.490338409425
$ws.Range("S6").Value = 0.4289341518184751
$ws.Range("T6").Value = 0.4289341518184751

$ws.Range("I7").Value = 0.7688940604785444
$ws.Range("J7").Value = 0.7688940604785441
$ws.Range("M7").Value = 36.057927
$ws.Range("N7").Value = 108.173781
$ws.Range("O7").Value = 0.1471761890142177
$ws.Range("P7").Value = 0.1471761890142177
$ws.Range("Q7").Value = 501.921608232254
$ws.Range("R7").Value = 4517.294474090286
$ws.Range("S7").Value = 0.1131628975768996
$ws.Range("T7").Value = 0.1131628975768995

$ws.Range("G8").Value = 3.7916
$ws.Range("H8").Value = 11.3748
$ws.Range("I8").Value = 0.209437228864931
$ws.Range("J8").Value = 0.209437228864931
$ws.Range("M8").Value = 72.266001
$ws.Range("N8").Value = 216.798003
$ws.Range("O8").Value = 0.2949652269937106
$ws.Range("P8").Value = 0.2949652269937106
$ws.Range("Q8").Value = 274.0037693916
$ws.Range("R8").Value = 2466.0339245244
$ws.Range("S8").Value = 0.06177669975307808
$ws.Range("T8").Value = 0.06177669975307808

$ws.Range("G9").Value = 3.7916
$ws.Range("H9").Value = 11.3748
$ws.Range("I9").Value = 0.209437228864931
$ws.Range("J9").Value = 0.209437228864931
$ws.Range("N9").Value = 410.023338
$ws.Range("O9").Value = 0.5578585839920717
$ws.Range("P9").Value = 0.5578585839920718
$ws.Range("Q9").Value = 518.2148294536
$ws.Range("R9").Value = 4663.933465082399
$ws.Range("S9").Value = 0.1168363559298138
$ws.Range("T9").Value = 0.1168363559298139

$ws.Range("G10").Value = 3.7916
$ws.Range("H10").Value = 11.3748
$ws.Range("I10").Value = 0.209437228864931
$ws.Range("J10").Value = 0.209437228864931
$ws.Range("M10").Value = 36.057927
$ws.Range("N10").Value = 108.173781
$ws.Range("O10").Value = 0.1471761890142177
$ws.Range("P10").Value = 0.1471761890142177
$ws.Range("Q10").Value = 136.7172360132
$ws.Range("R10").Value = 1230.4551241188
$ws.Range("S10").Value = 0.03082417318203905
$ws.Range("T10").Value = 0.03082417318203905

Write-Host "done"